# "#5: cash & deposit done"
# Rework the 存款 (deposit) sheet (sheet4) so that:
#   - row 1 becomes a proper header row (it used to be a stray duplicate of
#     row 2's data) spanning B1:M1 with the field names used elsewhere in
#     this workbook.
#   - every data row (2-18) gains the same trailing metadata columns
#     (G:property_category, H:category, I:date, J:legislator_name,
#      K:legislator_id, L:source_file, M:index) that the other sheets
#     already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Row 1: header labels -------------------------------------------------
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# --- Rows 2-18: append the shared metadata columns ------------------------
# (columns A-F already hold the correct bank/deposit_type/currency/owner/
#  total values and are left untouched)
for ($row = 2; $row -le 18; $row++) {
    $idx = $ws.Cells.Item($row, 1).Value2

    $ws.Cells.Item($row, 7).Value  = "deposit"
    $ws.Cells.Item($row, 8).Value  = "normal"
    $ws.Cells.Item($row, 9).Value  = "2011-11-17"
    $ws.Cells.Item($row, 10).Value = "李鴻鈞"
    $ws.Cells.Item($row, 11).Value = 898
    $ws.Cells.Item($row, 12).Value = "tmp6101"
    $ws.Cells.Item($row, 13).Value = $idx
}
